# Apply commit: "Preparation of databases for screening and analysis
# (standarize variable names, delete inadecuate values...)"
#
# - Rename the worksheet tab
# - Translate / standardize the column header names (row 1), which also
#   renames the matching Excel Table ("Tabla535") column headers
# - Update the active cell / selection
# - Resize the data columns to their new (narrower) widths now that the
#   English header labels are shorter than the original Spanish ones

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename worksheet tab
$ws.Name = "Root dmg (20 reps.)"

# Standardize / translate header row (also updates the Table535 column names)
$ws.Range("A1").Value = "Date"
$ws.Range("B1").Value = "Field"
$ws.Range("C1").Value = "Treatment"
$ws.Range("D1").Value = "Repeat"
$ws.Range("E1").Value = "Root_weight"
$ws.Range("F1").Value = "Observations"

# Resize columns now that headers are shorter (bestFit no longer needed)
$ws.Columns.Item(2).ColumnWidth = 8.619791666666666
$ws.Columns.Item(3).ColumnWidth = 12.166666666666666
$ws.Columns.Item(4).ColumnWidth = 9.256510416666666
$ws.Columns.Item(5).ColumnWidth = 16.619791666666668
$ws.Columns.Item(6).ColumnWidth = 93.07291666666667

# Update current selection / active cell
$ws.Range("H17").Select() | Out-Null
